$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.079.53"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.832.59"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'243.64"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'0.6277"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.07523"
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").Value = "'0.2923"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'23.17"
$ws.Range("E10").Value = "  +2.66%  "
$ws.Range("D11").Value = "'0.07720"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "1.830.67"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "'5.000"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").Value = "'0.6674"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "'82.57"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "'0.000009409"
$ws.Range("E16").Value = "  -7.16%  "
$ws.Range("D17").Value = "'5.982"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "29.092.29"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "2.080.62"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").Value = "'223.56"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "'7.110"
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'159.77"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").Value = "'0.1396"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").Value = "'8.493"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'17.93"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "'1.502"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").Value = "'0.05670"
$ws.Range("E30").Value = "  +8.75%  "
$ws.Range("D31").Value = "'4.152"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").Value = "'4.061"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").Value = "'1.202"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").Value = "'0.7486"
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("D35").Value = "'1.843"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'1.137"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "'2.667"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").Value = "'2.762"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "1.220.41"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").Value = "'0.01781"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("D41").Value = "'6.509"
$ws.Range("E41").Value = "  +2.73%  "
$ws.Range("D42").Value = "'0.8902"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "'102.03"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").Value = "1.979.90"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'65.51"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("D47").Value = "'0.00000000122"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").Value = "'0.5091"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").Value = "'0.07609"
$ws.Range("E49").Value = "  +11.92%  "
$ws.Range("D50").Value = "'0.4073"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").Value = "'9.021"
$ws.Range("E51").Value = "  +1.92%  "
